$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 55
$ws.Range("F4").Value = 20
$ws.Range("F5").Value = 211
$ws.Range("F6").Value = 1149
$ws.Range("G6").Value = 60
$ws.Range("F8").Value = 8437
$ws.Range("F11").Value = 6979
$ws.Range("F12").Value = 179
$ws.Range("F13").Value = 307
$ws.Range("F14").Value = 5149
$ws.Range("F15").Value = 5149
$ws.Range("F18").Value = 5671
$ws.Range("F19").Value = 5671
$ws.Range("F21").Value = 349
$ws.Range("F22").Value = 361
$ws.Range("F24").Value = 496
$ws.Range("F26").Value = 260
$ws.Range("F27").Value = 138
$ws.Range("F28").Value = 9467
$ws.Range("F29").Value = 77
$ws.Range("F30").Value = 1743
$ws.Range("F31").Value = 1185
$ws.Range("F32").Value = 40
$ws.Range("F33").Value = 1929
$ws.Range("F37").Value = 1019
$ws.Range("F38").Value = 1939
$ws.Range("F39").Value = 250
$ws.Range("F40").Value = 1251
$ws.Range("F41").Value = 58
$ws.Range("F42").Value = 4937
$ws.Range("F45").Value = 546
$ws.Range("F49").Value = 935
$ws.Range("F50").Value = 1293

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 5
$ws.Range("F11").Value = 185
$ws.Range("F19").Value = 897
$ws.Range("F22").Value = 2

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 55
$ws.Range("F5").Value = 20
$ws.Range("F6").Value = 211
$ws.Range("F8").Value = 1149
$ws.Range("G8").Value = 60
$ws.Range("F10").Value = 8437
$ws.Range("F13").Value = 6979
$ws.Range("F14").Value = 179
$ws.Range("F15").Value = 307
$ws.Range("F18").Value = 5149
$ws.Range("F19").Value = 5149
$ws.Range("F21").Value = 5671
$ws.Range("F22").Value = 5671
$ws.Range("F24").Value = 349
$ws.Range("F25").Value = 361
$ws.Range("F26").Value = 496
$ws.Range("F28").Value = 260
$ws.Range("F29").Value = 138
$ws.Range("F30").Value = 185
$ws.Range("F31").Value = 9468
$ws.Range("F32").Value = 77
$ws.Range("F33").Value = 1743
$ws.Range("F34").Value = 1185
$ws.Range("F35").Value = 40
$ws.Range("F36").Value = 1929
$ws.Range("F39").Value = 1019
$ws.Range("F40").Value = 1939
$ws.Range("F41").Value = 250
$ws.Range("F42").Value = 1251
$ws.Range("F43").Value = 4937
$ws.Range("F46").Value = 546
$ws.Range("F50").Value = 935
$ws.Range("F51").Value = 1293
